# New Commit Added On TimeStamp 01/06/2025 21:59
# Update the UploadFilePath values on the PracticeFormData sheet (rows 2-6, column N)
# to reflect the newly captured validation screenshot filenames.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PracticeFormData")

$ws.Range("N2").Value = "C:\Users\HP\Documents\Ahmedabad\ANKUSH_FRAMEWORK\ScreenShots\Validation\valdn01_10_2025_09_10_21.png"
$ws.Range("N3").Value = "C:\Users\HP\Documents\Ahmedabad\ANKUSH_FRAMEWORK\ScreenShots\Validation\valdn01_10_2025_09_10_59.png"
$ws.Range("N4").Value = "C:\Users\HP\Documents\Ahmedabad\ANKUSH_FRAMEWORK\ScreenShots\Validation\valdn01_12_2025_09_12_12.png"
$ws.Range("N5").Value = "C:\Users\HP\Documents\Ahmedabad\ANKUSH_FRAMEWORK\ScreenShots\Validation\valdn01_12_2025_09_12_45.png"
$ws.Range("N6").Value = "C:\Users\HP\Documents\Ahmedabad\ANKUSH_FRAMEWORK\ScreenShots\Validation\valdn01_13_2025_09_13_26.png"

$wb.Save()
